# STS IR Bot Performer / Data / Config.xlsx
#
# Commit: "Delete old files. Fixes after the demo with Jaime.
#          Click on Return was improved. Transaction reference:
#          Legal entity was added"
#
# The only functional/data change in the diff is on the "Constants" sheet:
# a brand-new configuration entry (key/value pair) is inserted right
# above the old row 49 ("StateBalancing_ReturnNamesPossibilities"),
# which pushes every following row down by one:
#
#   A49 = OtherDeductionsList_NonLocalStates
#   B49 = Non-Local State Forms
#
# This sits naturally right after the existing
# "OtherDeductionsList_SheetJeffersonCountyExceptionalCase" row (48),
# i.e. another OtherDeductionsList_* setting — matching the "Transaction
# reference: Legal entity was added" note in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert a new blank row above row 49; this shifts the previous rows
# 49..758 down to 50..759 (and extends the trailing blank rows, same as
# the real edit did).
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row with the new setting name/value pair.
$ws.Range("A49").Value = "OtherDeductionsList_NonLocalStates"
$ws.Range("B49").Value = "Non-Local State Forms"

# Match the row height used by every other data row on this sheet
# (14.25pt custom height) for both the new row and the row it pushed down.
$ws.Rows.Item(49).RowHeight = 14.25
$ws.Rows.Item(50).RowHeight = 14.25

# Reflect the updated selection/view saved in the workbook: the author's
# cursor ended up on the freshly-added cell A49.
$ws.Activate() | Out-Null
$ws.Range("A49").Select() | Out-Null
